# Drivers.xlsx — re-sort driver list by Driver ID (column A) and add a
# stray quote-prefixed placeholder value in E22, matching the author's
# "sort + scratch note" edit captured in the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Re-sort the data block (rows 2-27, cols A:D) ascending by
#        Driver ID (column A) — mirrors Data > Sort A-Z on column A.
#        Using the Sort object (rather than Range.Sort) so the sheet's
#        persisted sortState/sortCondition reflect the new range.
$ws.Sort.SortFields.Clear() | Out-Null
$ws.Sort.SortFields.Add($ws.Range("A1"), 0, 1, 0, 0) | Out-Null
$ws.Sort.SetRange($ws.Range("A2:D29")) | Out-Null
$ws.Sort.Header = 0
$ws.Sort.Apply() | Out-Null

# --- 2) Add the new E22 cell (FABIO LEIMER's row once re-sorted).
#        Copy the formatting from D22 first so the new cell inherits the
#        existing centered / wrapped style, then type a value that starts
#        with an apostrophe — Excel stores this as text with a
#        quote-prefix style (a new cellXfs entry gets created for it).
$ws.Range("D22").Copy() | Out-Null
$ws.Range("E22").PasteSpecial(-4122) | Out-Null
# Note: the leading apostrophe is consumed as Excel's quote-prefix marker
# (forces text, not stored), so 37 apostrophes here yields 36 stored.
$ws.Range("E22").Value = "'''''''''''''''''''''''''''''''''''''"

# --- 3) Update the view: scroll so column C is left-most and select E23,
#        matching where the author's cursor ended up after the edit.
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("E23").Select() | Out-Null
